$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H68").Value = 9900
$ws.Range("I68").Value = 9900
$ws.Range("K68").Value = 9900
$ws.Range("M68").Value = -9151
$ws.Range("H71").Value = 9900
$ws.Range("I71").Value = 9900
$ws.Range("K71").Value = 29700
$ws.Range("M71").Value = -25956
$ws.Range("H76").Value = 3370560
$ws.Range("I76").Value = 3707060.5
$ws.Range("K76").Value = 3707060.5
$ws.Range("M76").Value = -3706745.5
$ws.Range("H79").Value = 3370560
$ws.Range("I79").Value = 3707060.5
$ws.Range("K79").Value = 3707060.5
$ws.Range("M79").Value = -3705968.5
$ws.Range("H121").Value = 657.1667
$ws.Range("J121").Value = 607.4783
$ws.Range("L121").Value = 1822.4349
$ws.Range("N121").Value = -5316.4349
$ws.Range("H137").Value = 1553.25
$ws.Range("I137").Value = 1456.625
$ws.Range("J137").Value = 1649.875
$ws.Range("K137").Value = 4369.875
$ws.Range("L137").Value = 4949.625
$ws.Range("M137").Value = -1819.875
$ws.Range("N137").Value = -10049.625
$ws.Range("H138").Value = 2552.5693
$ws.Range("I138").Value = 1704.5758
$ws.Range("J138").Value = 3427.0625
$ws.Range("K138").Value = 5113.7274
$ws.Range("L138").Value = 10281.1875
$ws.Range("M138").Value = 26.27260000000024
$ws.Range("N138").Value = -20561.1875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 393601.5
$ws.Range("I32").Value = 5036.8066
$ws.Range("K32").Value = 5036.8066
$ws.Range("M32").Value = -4749.8066
$ws.Range("H45").Value = 2552.4546
$ws.Range("I45").Value = 2358.2144
$ws.Range("J45").Value = 2892.375
$ws.Range("K45").Value = 2358.2144
$ws.Range("L45").Value = 2892.375
$ws.Range("M45").Value = -1981.2144
$ws.Range("N45").Value = -3646.375
$ws.Range("H61").Value = 5294.7617
$ws.Range("I61").Value = 5431.0527
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 5431.0527
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -5219.0527
$ws.Range("N61").Value = -4424
$ws.Range("H74").Value = 1240.5385
$ws.Range("I74").Value = 818.1429000000001
$ws.Range("J74").Value = 1733.3334
$ws.Range("K74").Value = 818.1429000000001
$ws.Range("L74").Value = 1733.3334
$ws.Range("M74").Value = 55.85709999999995
$ws.Range("N74").Value = -3481.3334
$ws.Range("H77").Value = 1240.5385
$ws.Range("I77").Value = 818.1429000000001
$ws.Range("J77").Value = 1733.3334
$ws.Range("K77").Value = 4090.7145
$ws.Range("L77").Value = 8666.666999999999
$ws.Range("M77").Value = 277.2855
$ws.Range("N77").Value = -17402.667
$ws.Range("H88").Value = 5313.278
$ws.Range("I88").Value = 1999
$ws.Range("J88").Value = 5508.2354
$ws.Range("K88").Value = 1999
$ws.Range("L88").Value = 5508.2354
$ws.Range("M88").Value = -1593
$ws.Range("N88").Value = -6320.2354
$ws.Range("H91").Value = 5313.278
$ws.Range("I91").Value = 1999
$ws.Range("J91").Value = 5508.2354
$ws.Range("K91").Value = 1999
$ws.Range("L91").Value = 5508.2354
$ws.Range("M91").Value = -595
$ws.Range("N91").Value = -8316.2354
$ws.Range("H122").Value = 33929.477
$ws.Range("I122").Value = 38907.055
$ws.Range("J122").Value = 4064
$ws.Range("K122").Value = 116721.165
$ws.Range("L122").Value = 12192
$ws.Range("M122").Value = -114271.165
$ws.Range("N122").Value = -17092
$ws.Range("H132").Value = 1683
$ws.Range("I132").Value = 1409.2322
$ws.Range("J132").Value = 3076.7273
$ws.Range("K132").Value = 4227.696599999999
$ws.Range("L132").Value = 9230.1819
$ws.Range("M132").Value = -1697.696599999999
$ws.Range("N132").Value = -14290.1819
$ws.Range("H136").Value = 5294.7617
$ws.Range("I136").Value = 5431.0527
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 16293.1581
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -13743.1581
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 12794.161
$ws.Range("J86").Value = 3535.35
$ws.Range("L86").Value = 3535.35
$ws.Range("N86").Value = -5781.35
$ws.Range("H89").Value = 12794.161
$ws.Range("J89").Value = 3535.35
$ws.Range("L89").Value = 17676.75
$ws.Range("N89").Value = -28908.75
$ws.Range("H105").Value = 1519.4706
$ws.Range("I105").Value = 1479.8
$ws.Range("J105").Value = 1576.1428
$ws.Range("K105").Value = 1479.8
$ws.Range("L105").Value = 1576.1428
$ws.Range("M105").Value = 267.2
$ws.Range("N105").Value = -5070.1428
$ws.Range("H107").Value = 3473345
$ws.Range("I107").Value = 4065987.8
$ws.Range("K107").Value = 4065987.8
$ws.Range("M107").Value = -4064067.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2138.9
$ws.Range("I58").Value = 2900
$ws.Range("J58").Value = 2054.3333
$ws.Range("K58").Value = 2900
$ws.Range("L58").Value = 2054.3333
$ws.Range("M58").Value = -2697
$ws.Range("N58").Value = -2460.3333
$ws.Range("H99").Value = 2382.818
$ws.Range("I99").Value = 2153
$ws.Range("J99").Value = 2995.6667
$ws.Range("K99").Value = 2153
$ws.Range("L99").Value = 2995.6667
$ws.Range("M99").Value = -655
$ws.Range("N99").Value = -5991.6667
$ws.Range("H126").Value = 2382.818
$ws.Range("I126").Value = 2153
$ws.Range("J126").Value = 2995.6667
$ws.Range("K126").Value = 6459
$ws.Range("L126").Value = 8987.000100000001
$ws.Range("M126").Value = -3989
$ws.Range("N126").Value = -13927.0001
$ws.Range("H136").Value = 2138.9
$ws.Range("I136").Value = 2900
$ws.Range("J136").Value = 2054.3333
$ws.Range("K136").Value = 8700
$ws.Range("L136").Value = 6162.999899999999
$ws.Range("M136").Value = -6150
$ws.Range("N136").Value = -11262.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1666.6666
$ws.Range("I68").Value = 1000
$ws.Range("J68").Value = 2200
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 6600
$ws.Range("M68").Value = -2189
$ws.Range("N68").Value = -8222
$ws.Range("H71").Value = 1666.6666
$ws.Range("I71").Value = 1000
$ws.Range("J71").Value = 2200
$ws.Range("K71").Value = 9000
$ws.Range("L71").Value = 19800
$ws.Range("M71").Value = -4944
$ws.Range("N71").Value = -27912
$ws.Range("H113").Value = 979.0106
$ws.Range("I113").Value = 655.55554
$ws.Range("J113").Value = 1013.25885
$ws.Range("K113").Value = 1966.66662
$ws.Range("L113").Value = 3039.77655
$ws.Range("M113").Value = 203.33338
$ws.Range("N113").Value = -7379.77655
$ws.Range("H122").Value = 360.16
$ws.Range("I122").Value = 360.16
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3241.44
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -791.4400000000001
$ws.Range("N122").ClearContents()
$ws.Range("H131").Value = 5682713.5
$ws.Range("I131").Value = 1277.5
$ws.Range("J131").Value = 6579782
$ws.Range("K131").Value = 3832.5
$ws.Range("L131").Value = 19739346
$ws.Range("M131").Value = 1207.5
$ws.Range("N131").Value = -19749426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2138.7556
$ws.Range("I132").Value = 1614.3235
$ws.Range("J132").Value = 3759.7273
$ws.Range("K132").Value = 4842.970499999999
$ws.Range("L132").Value = 11279.1819
$ws.Range("M132").Value = -2312.970499999999
$ws.Range("N132").Value = -16339.1819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 7000
$ws.Range("J50").Value = 7000
$ws.Range("L50").Value = 7000
$ws.Range("N50").Value = -8274

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 200002400
$ws.Range("J14").Value = 2995
$ws.Range("L14").Value = 2995
$ws.Range("N14").Value = -3331
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
